# This edit relocates the bold "Play Arcade Bomb Free..." line and the
# italic review blurb that originally sat at the very end of the document:
#   * the bold line moves up to become a new paragraph right after the H1
#     title, its text becomes "Meta description" and the old blurb text is
#     appended after it (in a plain, non-bold run) as ": <blurb>";
#   * the italic paragraph that used to follow the bold line (now the very
#     last paragraph in the document) gets its text replaced with a new
#     image-generation prompt, keeping the italic run formatting.

$d = $word.ActiveDocument

$titleText  = "Play Arcade Bomb Free - Classic Slot Game without Complications"
$blurbText  = "Read our review of Arcade Bomb, a classic slot game with simple features and explosive bombs. Play for free and enjoy the fruit machine design."
$promptText = 'Prompt: Create a cartoon-style feature image that showcases the explosive fun of "Arcade Bomb". The image should feature a happy Maya warrior with glasses, who is surrounded by exploding bombs and fruits. Make the image lively and fun, with bright colors that pop. The Maya warrior should be the centerpiece of the image, with a beaming smile that shows how much fun he is having playing the game. In the background, include a space-themed backdrop with lots of stars and planets. The overall vibe of the image should be vibrant and exciting, capturing the essence of "Arcade Bomb".'

# --- locate the bold, non-heading "Play Arcade Bomb Free..." paragraph ---
$boldParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text.TrimEnd([char]13)
    if ($t -eq $titleText -and $para.Style.NameLocal -ne "Heading 1") {
        $boldParaIndex = $i
    }
}

$boldPara = $d.Paragraphs($boldParaIndex)
$boldRange = $boldPara.Range
$boldRange.Cut()

# --- insert a new (Normal-styled) paragraph right after the H1 title and
#     move the cut content -- together with its run/formatting structure --
#     into it ---
$titleRange = $d.Paragraphs.First.Range
$titleRange.Collapse(0)
$titleRange.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaPara.Style = "Normal"
$metaRange = $d.Paragraphs(2).Range
$metaRange.Paste()

# --- retarget the (still bold) run's text from the title to "Meta description" ---
$metaRange2 = $d.Paragraphs(2).Range
$boldTextRange = $d.Range($metaRange2.Start, $metaRange2.Start + $titleText.Length)
$boldTextRange.Text = "Meta description"

# --- append the (non-bold) blurb sentence right after it, before the
#     paragraph mark ---
$metaRange3 = $d.Paragraphs(2).Range
$insertPos = $metaRange3.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter(": " + $blurbText)

# --- replace the italic paragraph's text (now the last paragraph in the
#     document) with the new image-generation prompt, preserving italics ---
$blurbParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text.TrimEnd([char]13)
    if ($t -eq $blurbText) {
        $blurbParaIndex = $i
    }
}

$blurbPara = $d.Paragraphs($blurbParaIndex)
$br = $blurbPara.Range
$textRange = $d.Range($br.Start, $br.End - 1)
$textRange.Text = $promptText
